# Natmi following Dr Hou advice:
# Recompute the Pcsk9-Vldlr ligand-receptor edge table across all three
# cell-type clusters (FAPs, sCs, and the newly-added ECs), expanding the
# sheet from a 2x2 (4-row) grid to a full 3x3 (9-row) grid of
# sending-cluster x target-cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Pcsk9"
$ws.Cells.Item(2, 3).Value = "Vldlr"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.286724
$ws.Cells.Item(2, 8).Value = 0.8601719999999999
$ws.Cells.Item(2, 9).Value = 0.1240995385141648
$ws.Cells.Item(2, 10).Value = 0.1240995385141648
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.3374003333333333
$ws.Cells.Item(2, 14).Value = 1.012201
$ws.Cells.Item(2, 15).Value = 0.01738364872808817
$ws.Cells.Item(2, 16).Value = 0.01738364872808818
$ws.Cells.Item(2, 17).Value = 0.09674077317466664
$ws.Cells.Item(2, 18).Value = 0.8706669585719998
$ws.Cells.Item(2, 19).Value = 0.00215730278484809
$ws.Cells.Item(2, 20).Value = 0.00215730278484809

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Pcsk9"
$ws.Cells.Item(3, 3).Value = "Vldlr"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.286724
$ws.Cells.Item(3, 8).Value = 0.8601719999999999
$ws.Cells.Item(3, 9).Value = 0.1240995385141648
$ws.Cells.Item(3, 10).Value = 0.1240995385141648
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 16.71131166666667
$ws.Cells.Item(3, 14).Value = 50.133935
$ws.Cells.Item(3, 15).Value = 0.8610055862391021
$ws.Cells.Item(3, 16).Value = 0.8610055862391023
$ws.Cells.Item(3, 17).Value = 4.791534126313333
$ws.Cells.Item(3, 18).Value = 43.12380713682
$ws.Cells.Item(3, 19).Value = 0.1068503959103905
$ws.Cells.Item(3, 20).Value = 0.1068503959103905

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Pcsk9"
$ws.Cells.Item(4, 3).Value = "Vldlr"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.286724
$ws.Cells.Item(4, 8).Value = 0.8601719999999999
$ws.Cells.Item(4, 9).Value = 0.1240995385141648
$ws.Cells.Item(4, 10).Value = 0.1240995385141648
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.360351
$ws.Cells.Item(4, 14).Value = 7.081053000000001
$ws.Cells.Item(4, 15).Value = 0.1216107650328097
$ws.Cells.Item(4, 16).Value = 0.1216107650328097
$ws.Cells.Item(4, 17).Value = 0.676769280124
$ws.Cells.Item(4, 18).Value = 6.090923521116
$ws.Cells.Item(4, 19).Value = 0.0150918398189262
$ws.Cells.Item(4, 20).Value = 0.01509183981892621

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Pcsk9"
$ws.Cells.Item(5, 3).Value = "Vldlr"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.073738333333333
$ws.Cells.Item(5, 8).Value = 3.221215
$ws.Cells.Item(5, 9).Value = 0.4647341403288008
$ws.Cells.Item(5, 10).Value = 0.4647341403288009
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.3374003333333333
$ws.Cells.Item(5, 14).Value = 1.012201
$ws.Cells.Item(5, 15).Value = 0.01738364872808817
$ws.Cells.Item(5, 16).Value = 0.01738364872808818
$ws.Cells.Item(5, 17).Value = 0.3622796715794444
$ws.Cells.Item(5, 18).Value = 3.260517044215
$ws.Cells.Item(5, 19).Value = 0.00807877504742591
$ws.Cells.Item(5, 20).Value = 0.008078775047425914

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Pcsk9"
$ws.Cells.Item(6, 3).Value = "Vldlr"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.073738333333333
$ws.Cells.Item(6, 8).Value = 3.221215
$ws.Cells.Item(6, 9).Value = 0.4647341403288008
$ws.Cells.Item(6, 10).Value = 0.4647341403288009
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 16.71131166666667
$ws.Cells.Item(6, 14).Value = 50.133935
$ws.Cells.Item(6, 15).Value = 0.8610055862391021
$ws.Cells.Item(6, 16).Value = 0.8610055862391023
$ws.Cells.Item(6, 17).Value = 17.94357593678055
$ws.Cells.Item(6, 18).Value = 161.492183431025
$ws.Cells.Item(6, 19).Value = 0.4001386909391243
$ws.Cells.Item(6, 20).Value = 0.4001386909391245

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Pcsk9"
$ws.Cells.Item(7, 3).Value = "Vldlr"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.073738333333333
$ws.Cells.Item(7, 8).Value = 3.221215
$ws.Cells.Item(7, 9).Value = 0.4647341403288008
$ws.Cells.Item(7, 10).Value = 0.4647341403288009
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.360351
$ws.Cells.Item(7, 14).Value = 7.081053000000001
$ws.Cells.Item(7, 15).Value = 0.1216107650328097
$ws.Cells.Item(7, 16).Value = 0.1216107650328097
$ws.Cells.Item(7, 17).Value = 2.534399348821667
$ws.Cells.Item(7, 18).Value = 22.809594139395
$ws.Cells.Item(7, 19).Value = 0.05651667434225058
$ws.Cells.Item(7, 20).Value = 0.05651667434225062

# Row 8
$ws.Cells.Item(8, 1).Value = "ECs"
$ws.Cells.Item(8, 2).Value = "Pcsk9"
$ws.Cells.Item(8, 3).Value = "Vldlr"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9499733333333333
$ws.Cells.Item(8, 8).Value = 2.84992
$ws.Cells.Item(8, 9).Value = 0.4111663211570343
$ws.Cells.Item(8, 10).Value = 0.4111663211570343
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.3374003333333333
$ws.Cells.Item(8, 14).Value = 1.012201
$ws.Cells.Item(8, 15).Value = 0.01738364872808817
$ws.Cells.Item(8, 16).Value = 0.01738364872808818
$ws.Cells.Item(8, 17).Value = 0.3205213193244444
$ws.Cells.Item(8, 18).Value = 2.88469187392
$ws.Cells.Item(8, 19).Value = 0.007147570895814173
$ws.Cells.Item(8, 20).Value = 0.007147570895814175

# Row 9
$ws.Cells.Item(9, 1).Value = "ECs"
$ws.Cells.Item(9, 2).Value = "Pcsk9"
$ws.Cells.Item(9, 3).Value = "Vldlr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9499733333333333
$ws.Cells.Item(9, 8).Value = 2.84992
$ws.Cells.Item(9, 9).Value = 0.4111663211570343
$ws.Cells.Item(9, 10).Value = 0.4111663211570343
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 16.71131166666667
$ws.Cells.Item(9, 14).Value = 50.133935
$ws.Cells.Item(9, 15).Value = 0.8610055862391021
$ws.Cells.Item(9, 16).Value = 0.8610055862391023
$ws.Cells.Item(9, 17).Value = 15.87530044835556
$ws.Cells.Item(9, 18).Value = 142.8777040352
$ws.Cells.Item(9, 19).Value = 0.3540164993895872
$ws.Cells.Item(9, 20).Value = 0.3540164993895873

# Row 10
$ws.Cells.Item(10, 1).Value = "ECs"
$ws.Cells.Item(10, 2).Value = "Pcsk9"
$ws.Cells.Item(10, 3).Value = "Vldlr"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9499733333333333
$ws.Cells.Item(10, 8).Value = 2.84992
$ws.Cells.Item(10, 9).Value = 0.4111663211570343
$ws.Cells.Item(10, 10).Value = 0.4111663211570343
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.360351
$ws.Cells.Item(10, 14).Value = 7.081053000000001
$ws.Cells.Item(10, 15).Value = 0.1216107650328097
$ws.Cells.Item(10, 16).Value = 0.1216107650328097
$ws.Cells.Item(10, 17).Value = 2.242270507306667
$ws.Cells.Item(10, 18).Value = 20.18043456576
$ws.Cells.Item(10, 19).Value = 0.05000225087163285
$ws.Cells.Item(10, 20).Value = 0.05000225087163288
